# LOS Galacticos roster update.
#
# The underlying (player, position, team) rows get reshuffled: several
# players move to different rows while keeping their own position/team
# data intact. Net effect, expressed as per-cell writes against the
# original layout:
#   - Row 2  becomes Jose Alvarado        / PG       / New Orleans Pelicans
#   - Row 7  becomes Giannis Antetokounmpo/ PF,C     / Milwaukee Bucks
#   - Row 9  becomes Jaren Jackson Jr.    / (PF,C)   / Memphis Grizzlies
#   - Row 11 becomes Draymond Green       / (PF,C)   / Golden State Warriors
#   - Row 14 becomes Keyonte George       / PG,SG    / Utah Jazz
#   - Row 15 becomes Amen Thompson        / SG,SF,PF / Houston Rockets
#   - Row 16 becomes Paul George          / (SG,SF,PF)/ Philadelphia 76ers
# Rows 3,4,5,6,8,10,12,13,17,18 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jose Alvarado"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "New Orleans Pelicans"

$ws.Range("A7").Value = "Giannis Antetokounmpo"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Milwaukee Bucks"

$ws.Range("A9").Value = "Jaren Jackson Jr."
$ws.Range("C9").Value = "Memphis Grizzlies"

$ws.Range("A11").Value = "Draymond Green"
$ws.Range("C11").Value = "Golden State Warriors"

$ws.Range("A14").Value = "Keyonte George"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Utah Jazz"

$ws.Range("A15").Value = "Amen Thompson"
$ws.Range("B15").Value = "SG,SF,PF"
$ws.Range("C15").Value = "Houston Rockets"

$ws.Range("A16").Value = "Paul George"
$ws.Range("C16").Value = "Philadelphia 76ers"
